$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns hold plain text in the source
# workbook (no numeric/percentage cell formatting). Temporarily force the
# target range to Text format so the refreshed values are written back as
# text rather than being auto-coerced into numbers/percentages, then strip
# the temporary formatting again so the cells end up styled exactly like
# they started (default/general style).
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "289.61"
$ws.Range("E2").Value = "-9.18%"
$ws.Range("D3").Value = "40.40"
$ws.Range("E3").Value = "-2.31%"
$ws.Range("D4").Value = "5.048"
$ws.Range("E4").Value = "-3.66%"
$ws.Range("D5").Value = "0.07287"
$ws.Range("E5").Value = "-5.44%"
$ws.Range("D6").Value = "4.282"
$ws.Range("E6").Value = "-1.42%"
$ws.Range("D7").Value = "1.561"
$ws.Range("E7").Value = "-7.16%"
$ws.Range("D8").Value = "0.9205"
$ws.Range("E8").Value = "-2.40%"
$ws.Range("D9").Value = "0.1155"
$ws.Range("E9").Value = "-8.37%"
$ws.Range("D10").Value = "0.1723"
$ws.Range("E10").Value = "-7.56%"
$ws.Range("D11").Value = "0.08672"
$ws.Range("E11").Value = "-5.42%"
$ws.Range("D12").Value = "0.04180"
$ws.Range("E12").Value = "1.45%"
$ws.Range("E13").Value = "0.28%"
$ws.Range("D14").Value = "0.001272"
$ws.Range("E14").Value = "-1.12%"
$ws.Range("D15").Value = "0.005902"
$ws.Range("E15").Value = "-2.31%"
$ws.Range("D16").Value = "3.396"
$ws.Range("E16").Value = "1.41%"
$ws.Range("D19").Value = "7.883"
$ws.Range("E19").Value = "-6.32%"
$ws.Range("D20").Value = "0.1351"
$ws.Range("E20").Value = "-0.27%"
$ws.Range("D21").Value = "0.2884"
$ws.Range("E21").Value = "5.55%"
$ws.Range("D22").Value = "0.03872"
$ws.Range("E22").Value = "-4.24%"
$ws.Range("D23").Value = "0.001270"
$ws.Range("E23").Value = "-0.17%"
$ws.Range("D24").Value = "0.003845"
$ws.Range("E24").Value = "-6.67%"
$ws.Range("D25").Value = "0.0001282"
$ws.Range("E25").Value = "0.55%"
$ws.Range("D26").Value = "0.0003728"
$ws.Range("E26").Value = "-95.02%"
$ws.Range("D38").Value = "0.02317"
$ws.Range("E38").Value = "-8.64%"
$ws.Range("D39").Value = "0.04973"
$ws.Range("E39").Value = "-6.49%"
$ws.Range("D40").Value = "0.006641"
$ws.Range("E40").Value = "208.16%"
$ws.Range("D41").Value = "0.007679"
$ws.Range("E41").Value = "-1.35%"
$ws.Range("E42").Value = "-3.19%"
$ws.Range("D43").Value = "0.007372"
$ws.Range("E43").Value = "4.71%"
$ws.Range("D44").Value = "0.007066"
$ws.Range("E44").Value = "-14.90%"
$ws.Range("D45").Value = "0.2893"
$ws.Range("E45").Value = "-16.92%"
$ws.Range("D46").Value = "0.00006416"
$ws.Range("E46").Value = "-3.99%"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("E48").Value = "-90.25%"
$ws.Range("E49").Value = "-0.19%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "-0.07%"

$rng.ClearFormats()
